$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback datetimes for row 2
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-08 07:22:27"
$wsZh.Range("G2").Value = "2016-01-08 07:23:11"

# de-de sheet: update Correspond Handoff/Handback datetimes for row 2
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-08 07:22:40"
$wsDe.Range("G2").Value = "2016-01-08 07:23:29"
